$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 102,3
$data[0,0] = 63
$data[0,1] = '$(\eta_{p})_{3}$'
$data[0,2] = [double]"0.0004147835734525503"
$data[1,0] = 12
$data[1,1] = 'h$_{q}$'
$data[1,2] = [double]"0.0002908560973987411"
$data[2,0] = 34
$data[2,1] = '$F_{q}$'
$data[2,2] = [double]"0.0002254040881591021"
$data[3,0] = 71
$data[3,1] = '$\langle qq \vert qq \rangle$'
$data[3,2] = [double]"0.0001621271885981359"
$data[4,0] = 23
$data[4,1] = '$type_0$'
$data[4,2] = [double]"0.0001340690978451755"
$data[5,0] = 35
$data[5,1] = '$\eta_{q}$'
$data[5,2] = [double]"0.0001234474159007336"
$data[6,0] = 25
$data[6,1] = '$type_2$'
$data[6,2] = [double]"0.0001176476931310785"
$data[7,0] = 22
$data[7,1] = 'h$_{s}$'
$data[7,2] = [double]"0.000109611698848675"
$data[8,0] = 51
$data[8,1] = '$(\eta_{r})_{1}$'
$data[8,2] = [double]"0.0001092924295970475"
$data[9,0] = 33
$data[9,1] = '$\omega_{q}$'
$data[9,2] = [double]"7.583337057740694e-05"
$data[10,0] = 82
$data[10,1] = '$(\langle pq \vert pq \rangle)_{1}$'
$data[10,2] = [double]"7.030886905455172e-05"
$data[11,0] = 97
$data[11,1] = '$(\langle rr \vert rr \rangle)_{3}$'
$data[11,2] = [double]"5.540110250131223e-05"
$data[12,0] = 57
$data[12,1] = '$(\omega_{r})_{2}$'
$data[12,2] = [double]"5.528857922024196e-05"
$data[13,0] = 40
$data[13,1] = '$F_{s}^{\text{SCF}}$'
$data[13,2] = [double]"5.3907237382785e-05"
$data[14,0] = 41
$data[14,1] = '$\omega_{s}$'
$data[14,2] = [double]"5.067261043112274e-05"
$data[15,0] = 56
$data[15,1] = '$(F_{r}^{\text{SCF}})_{2}$'
$data[15,2] = [double]"5.022533929130881e-05"
$data[16,0] = 43
$data[16,1] = '$\eta_{s}$'
$data[16,2] = [double]"5.010879139527134e-05"
$data[17,0] = 37
$data[17,1] = '$(\omega_{r})_{0}$'
$data[17,2] = [double]"4.894964012305951e-05"
$data[18,0] = 64
$data[18,1] = '$(F_{r}^{\text{SCF}})_{3}$'
$data[18,2] = [double]"4.854578504549892e-05"
$data[19,0] = 49
$data[19,1] = '$(\omega_{r})_{1}$'
$data[19,2] = [double]"4.613033188808599e-05"
$data[20,0] = 53
$data[20,1] = '$(\omega_{p})_{2}$'
$data[20,2] = [double]"4.578931506235848e-05"
$data[21,0] = 32
$data[21,1] = '$F_{q}^{\text{SCF}}$'
$data[21,2] = [double]"4.315669946392124e-05"
$data[22,0] = 13
$data[22,1] = 'h$_{qs}$'
$data[22,2] = [double]"4.272819666856654e-05"
$data[23,0] = 73
$data[23,1] = '$\langle ss \vert ss \rangle$'
$data[23,2] = [double]"3.974723130416576e-05"
$data[24,0] = 96
$data[24,1] = '$(\langle pp \vert pp \rangle)_{3}$'
$data[24,2] = [double]"3.738932487678512e-05"
$data[25,0] = 0
$data[25,1] = '(h$_{p}$)$_{0}$'
$data[25,2] = [double]"3.510612985424193e-05"
$data[26,0] = 42
$data[26,1] = '$F_{s}$'
$data[26,2] = [double]"3.382181412097486e-05"
$data[27,0] = 26
$data[27,1] = '$type_3$'
$data[27,2] = [double]"3.352202308686562e-05"
$data[28,0] = 16
$data[28,1] = '(h$_{r}$)$_{2}$'
$data[28,2] = [double]"3.269623720854638e-05"
$data[29,0] = 99
$data[29,1] = '$(\langle pq \vert qp \rangle)_{3}$'
$data[29,2] = [double]"3.144901541523325e-05"
$data[30,0] = 90
$data[30,1] = '$(\langle pq \vert pq \rangle)_{2}$'
$data[30,2] = [double]"3.066937712270076e-05"
$data[31,0] = 78
$data[31,1] = '$(\langle pq \vert rs \rangle)_{1}$'
$data[31,2] = [double]"2.972545353383824e-05"
$data[32,0] = 98
$data[32,1] = '$(\langle pq \vert pq \rangle)_{3}$'
$data[32,2] = [double]"2.924205366555794e-05"
$data[33,0] = 46
$data[33,1] = '$(F_{p})_{1}$'
$data[33,2] = [double]"2.880797824305791e-05"
$data[34,0] = 52
$data[34,1] = '$(F_{p}^{\text{SCF}})_{2}$'
$data[34,2] = [double]"2.877772974373062e-05"
$data[35,0] = 65
$data[35,1] = '$(\omega_{r})_{3}$'
$data[35,2] = [double]"2.790245283656804e-05"
$data[36,0] = 17
$data[36,1] = '(h$_{r}$)$_{3}$'
$data[36,2] = [double]"2.762249918575351e-05"
$data[37,0] = 70
$data[37,1] = '$(\langle pp \vert pp \rangle)_{0}$'
$data[37,2] = [double]"2.555189428686214e-05"
$data[38,0] = 59
$data[38,1] = '$(\eta_{r})_{2}$'
$data[38,2] = [double]"2.500322265052189e-05"
$data[39,0] = 62
$data[39,1] = '$(F_{p})_{3}$'
$data[39,2] = [double]"2.467730605834434e-05"
$data[40,0] = 3
$data[40,1] = '(h$_{p}$)$_{3}$'
$data[40,2] = [double]"2.370859936291518e-05"
$data[41,0] = 89
$data[41,1] = '$(\langle rr \vert rr \rangle)_{2}$'
$data[41,2] = [double]"2.24425735379371e-05"
$data[42,0] = 67
$data[42,1] = '$(\eta_{r})_{3}$'
$data[42,2] = [double]"2.194952573238563e-05"
$data[43,0] = 101
$data[43,1] = '$(\langle rs \vert sr \rangle)_{3}$'
$data[43,2] = [double]"2.166553436885997e-05"
$data[44,0] = 2
$data[44,1] = '(h$_{p}$)$_{2}$'
$data[44,2] = [double]"2.161995342224316e-05"
$data[45,0] = 44
$data[45,1] = '$(F_{p}^{\text{SCF}})_{1}$'
$data[45,2] = [double]"2.145205228392148e-05"
$data[46,0] = 66
$data[46,1] = '$(F_{r})_{3}$'
$data[46,2] = [double]"2.131186471675095e-05"
$data[47,0] = 1
$data[47,1] = '(h$_{p}$)$_{1}$'
$data[47,2] = [double]"2.127935391801965e-05"
$data[48,0] = 74
$data[48,1] = '$(\langle pq \vert pq \rangle)_{0}$'
$data[48,2] = [double]"2.083144159553002e-05"
$data[49,0] = 28
$data[49,1] = '$(F_{p}^{\text{SCF}})_{0}$'
$data[49,2] = [double]"2.048424288643049e-05"
$data[50,0] = 61
$data[50,1] = '$(\omega_{p})_{3}$'
$data[50,2] = [double]"2.019384158366014e-05"
$data[51,0] = 48
$data[51,1] = '$(F_{r}^{\text{SCF}})_{1}$'
$data[51,2] = [double]"1.86468955907999e-05"
$data[52,0] = 24
$data[52,1] = '$type_1$'
$data[52,2] = [double]"1.856037332654823e-05"
$data[53,0] = 75
$data[53,1] = '$(\langle pq \vert qp \rangle)_{0}$'
$data[53,2] = [double]"1.825175037009249e-05"
$data[54,0] = 85
$data[54,1] = '$(\langle rs \vert sr \rangle)_{1}$'
$data[54,2] = [double]"1.782761198872627e-05"
$data[55,0] = 21
$data[55,1] = '(h$_{rs}$)$_{3}$'
$data[55,2] = [double]"1.649151046863715e-05"
$data[56,0] = 45
$data[56,1] = '$(\omega_{p})_{1}$'
$data[56,2] = [double]"1.614419651622889e-05"
$data[57,0] = 58
$data[57,1] = '$(F_{r})_{2}$'
$data[57,2] = [double]"1.405739590167229e-05"
$data[58,0] = 15
$data[58,1] = '(h$_{r}$)$_{1}$'
$data[58,2] = [double]"1.398789803635251e-05"
$data[59,0] = 93
$data[59,1] = '$(\langle rs \vert sr \rangle)_{2}$'
$data[59,2] = [double]"1.39204201783251e-05"
$data[60,0] = 92
$data[60,1] = '$(\langle rs\vert rs \rangle)_{2}$'
$data[60,2] = [double]"1.383879591033324e-05"
$data[61,0] = 100
$data[61,1] = '$(\langle rs\vert rs \rangle)_{3}$'
$data[61,2] = [double]"1.377250676613535e-05"
$data[62,0] = 47
$data[62,1] = '$(\eta_{p})_{1}$'
$data[62,2] = [double]"1.370057122698882e-05"
$data[63,0] = 54
$data[63,1] = '$(F_{p})_{2}$'
$data[63,2] = [double]"1.364546250502321e-05"
$data[64,0] = 83
$data[64,1] = '$(\langle pq \vert qp \rangle)_{1}$'
$data[64,2] = [double]"1.342526736476277e-05"
$data[65,0] = 55
$data[65,1] = '$(\eta_{p})_{2}$'
$data[65,2] = [double]"1.341475258890273e-05"
$data[66,0] = 84
$data[66,1] = '$(\langle rs\vert rs \rangle)_{1}$'
$data[66,2] = [double]"1.298724461177479e-05"
$data[67,0] = 29
$data[67,1] = '$(\omega_{p})_{0}$'
$data[67,2] = [double]"1.257822143379857e-05"
$data[68,0] = 81
$data[68,1] = '$(\langle rr \vert rr \rangle)_{1}$'
$data[68,2] = [double]"1.208708811872289e-05"
$data[69,0] = 91
$data[69,1] = '$(\langle pq \vert qp \rangle)_{2}$'
$data[69,2] = [double]"1.206194002070736e-05"
$data[70,0] = 60
$data[70,1] = '$(F_{p}^{\text{SCF}})_{3}$'
$data[70,2] = [double]"1.14480821947826e-05"
$data[71,0] = 72
$data[71,1] = '$(\langle rr \vert rr \rangle)_{0}$'
$data[71,2] = [double]"1.054587333829056e-05"
$data[72,0] = 88
$data[72,1] = '$(\langle pp \vert pp \rangle)_{2}$'
$data[72,2] = [double]"1.046330578617555e-05"
$data[73,0] = 50
$data[73,1] = '$(F_{r})_{1}$'
$data[73,2] = [double]"9.736810460318345e-06"
$data[74,0] = 7
$data[74,1] = '(h$_{pq}$)$_{3}$'
$data[74,2] = [double]"9.64952475843561e-06"
$data[75,0] = 30
$data[75,1] = '$(F_{p})_{0}$'
$data[75,2] = [double]"9.511353650007364e-06"
$data[76,0] = 5
$data[76,1] = '(h$_{pq}$)$_{1}$'
$data[76,2] = [double]"9.181144093710802e-06"
$data[77,0] = 6
$data[77,1] = '(h$_{pq}$)$_{2}$'
$data[77,2] = [double]"8.746495189779597e-06"
$data[78,0] = 77
$data[78,1] = '$(\langle rs \vert sr \rangle)_{0}$'
$data[78,2] = [double]"8.669908306326509e-06"
$data[79,0] = 36
$data[79,1] = '$(F_{r}^{\text{SCF}})_{0}$'
$data[79,2] = [double]"8.60787769624548e-06"
$data[80,0] = 76
$data[80,1] = '$(\langle rs\vert rs \rangle)_{0}$'
$data[80,2] = [double]"8.57747991604444e-06"
$data[81,0] = 80
$data[81,1] = '$(\langle pp \vert pp \rangle)_{1}$'
$data[81,2] = [double]"7.58459754194243e-06"
$data[82,0] = 39
$data[82,1] = '$(\eta_{r})_{0}$'
$data[82,2] = [double]"7.424434096315851e-06"
$data[83,0] = 20
$data[83,1] = '(h$_{rs}$)$_{2}$'
$data[83,2] = [double]"6.614931994552312e-06"
$data[84,0] = 19
$data[84,1] = '(h$_{rs}$)$_{1}$'
$data[84,2] = [double]"6.446747283783063e-06"
$data[85,0] = 10
$data[85,1] = '(h$_{pr}$)$_{2}$'
$data[85,2] = [double]"6.103454588417964e-06"
$data[86,0] = 14
$data[86,1] = '(h$_{r}$)$_{0}$'
$data[86,2] = [double]"5.285468178272739e-06"
$data[87,0] = 8
$data[87,1] = '(h$_{pr}$)$_{0}$'
$data[87,2] = [double]"5.130619698858806e-06"
$data[88,0] = 18
$data[88,1] = '(h$_{rs}$)$_{0}$'
$data[88,2] = [double]"4.861943880478145e-06"
$data[89,0] = 86
$data[89,1] = '$(\langle pq \vert rs \rangle)_{2}$'
$data[89,2] = [double]"4.702387727635698e-06"
$data[90,0] = 9
$data[90,1] = '(h$_{pr}$)$_{1}$'
$data[90,2] = [double]"4.542776859953937e-06"
$data[91,0] = 4
$data[91,1] = '(h$_{pq}$)$_{0}$'
$data[91,2] = [double]"4.358596413610488e-06"
$data[92,0] = 94
$data[92,1] = '$(\langle pq \vert rs \rangle)_{3}$'
$data[92,2] = [double]"4.044607385059703e-06"
$data[93,0] = 11
$data[93,1] = '(h$_{pr}$)$_{3}$'
$data[93,2] = [double]"3.371106522828484e-06"
$data[94,0] = 38
$data[94,1] = '$(F_{r})_{0}$'
$data[94,2] = [double]"2.74040683384151e-06"
$data[95,0] = 68
$data[95,1] = '$(\langle pq \vert rs \rangle)_{0}$'
$data[95,2] = [double]"1.565196796589627e-06"
$data[96,0] = 31
$data[96,1] = '$(\eta_{p})_{0}$'
$data[96,2] = [double]"1.530204272353371e-06"
$data[97,0] = 87
$data[97,1] = '$(\langle pq \vert sr \rangle)_{2}$'
$data[97,2] = [double]"2.231152104564879e-07"
$data[98,0] = 69
$data[98,1] = '$(\langle pq \vert sr \rangle)_{0}$'
$data[98,2] = [double]"2.139017410912953e-07"
$data[99,0] = 79
$data[99,1] = '$(\langle pq \vert sr \rangle)_{1}$'
$data[99,2] = [double]"2.002322563337279e-07"
$data[100,0] = 95
$data[100,1] = '$(\langle pq \vert sr \rangle)_{3}$'
$data[100,2] = [double]"1.723540907865675e-07"
$data[101,0] = 27
$data[101,1] = '$\mathbf{b}$'
$data[101,2] = [double]"2.709181034506403e-08"

$ws.Range("A2:C103").Value = $data